$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'61.571.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.22%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.388.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.74%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'577.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.07%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'137.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.20%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'3.386.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.75%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.472"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.26%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'7.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.05%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.125"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.59%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.389"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.11%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'3.968.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.81%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("E14").Value = "'  +2.09%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'0.0000176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.45%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = "'WrappedEther"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'3.394.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.91%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = "'Avalanche"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'25.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.63%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'61.674.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.11%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'14.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.32%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = "'  +1.05%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'9.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.09%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'376.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.36%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.558"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.59%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'3.526.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.85%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.17%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.0000126"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +7.94%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'71.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.78%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  -0.27%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'7.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.78%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.09%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'0.160"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.99%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'8.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.33%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  +1.51%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = "'  +0.07%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'23.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.20%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'5.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.15%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'1.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.55%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'6.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.20%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'165.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.22%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.0781"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.39%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = "'FirstDigitalUSD"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.14%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "'Mantle"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.781"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.83%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "'Stacks"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'1.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +8.72%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").Value = "'ONDO"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'1.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.15%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'25.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +8.99%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'4.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.14%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'41.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.24%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'6.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.48%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'22.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.05%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'2.345.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +6.28%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.0261"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.01%  "
$ws.Range("E51").Style = "Normal"
